# Findings Tracker - remove the per-row "Business Exception" logging
# entries (SOXITBOT / Travis Lee) that used to live in column M
# ("ITRCA Member that filed or reviewed(if BOT found) final evidence")
# for the data rows 2-17. The shared-string table drops the now-unused
# "SOXITBOT" / "Travis Lee" entries automatically once every cell that
# referenced them is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindingsTracker")

# Clear the Business Exception / reviewer values in column M for the
# data rows, leaving the cell formatting (style) untouched.
$ws.Range("M2:M17").ClearContents()

# Match the workbook's on-disk selection/view state left behind by the
# editor: column M selected for the same row range.
$ws.Activate()
$ws.Range("M2:M17").Select()
